# CSI1101 A2: finished password protection, todo choosing password
# Append four new bibliography rows (in-text citation / full reference pairs)
# to the end of the existing table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (inText citation, full reference) pairs, entered in the same order the
# author typed them -- column A first, then column B -- for the first three
# new rows. The final row was filled in with the full reference pasted in
# column B before the short in-text citation was typed into column A.
$ws.Cells.Item(13, 1).Value = "(MITRE, 2014)"
$ws.Cells.Item(13, 2).Value = "MITRE. (2014). CVE - CVE-2014-1532. cve.mitre.org. Retrieved March 18, 2015, from http://cve.mitre.org/cgi-bin/cvename.cgi?name=CVE-2014-1532"

$ws.Cells.Item(14, 1).Value = "(SecurityFocus, 2010)"
$ws.Cells.Item(14, 2).Value = "SecurityFocus. (2010). Adobe Reader and Acrobat U3D Support Remote Code Execution Vulnerability. Retrieved March 11, 2015, from http://www.securityfocus.com/bid/37756/info"

$ws.Cells.Item(15, 1).Value = "(SecurityTracker, 2010)"
$ws.Cells.Item(15, 2).Value = "SecurityTracker. (2010). Adobe Acrobat and Adobe Reader Flaws Lets Remote Users Execute Arbitrary Code and Deny Service - SecurityTracker. SecurityTracker. Retrieved March 11, 2015, from http://www.securitytracker.com/id?1023446"

$ws.Cells.Item(16, 2).Value = "MITRE. (2006). CVE - CVE-2006-2198. cve.mitre.org. Retrieved March 11, 2015, from http://cve.mitre.org/cgi-bin/cvename.cgi?name=CVE-2006-2198"
$ws.Cells.Item(16, 1).Value = "(MITRE, 2006)"

$ws.Range("B34").Select()
